# Add a new day's trade sheet ("28-5-23 - 1") after the last existing
# "27-5-23 - 3" sheet, following the same template (headers, Notes /
# Confluences merged blocks, column widths, styles) that every other
# day-sheet in this workbook already uses, then fill in the new trade
# row and attach its 4-hour chart screenshot.

$wb = $excel.ActiveWorkbook

# Clone the last day-sheet so the new one inherits its formatting,
# merged cells ("Notes:" / "Confluences:") and column widths exactly.
$template = $wb.Worksheets.Item("27-5-23 - 3")
$template.Copy($null, $template)
$ws = $wb.ActiveSheet
$ws.Name = "28-5-23 - 1"

# Clear out the template's sample trade row/picture reference and
# replace it with the new day's single AUD/USD trade.
$ws.Range("A2").Value = "AUD/USD"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = "buy"
$ws.Range("D2").Value = 0.65195
$ws.Range("E2").Value = 0.65195
$ws.Range("F2").Value = "£55.33"
$ws.Range("G2").Value = 26.36
$ws.Range("H2").Value = "£0"
$ws.Range("I2").Value = "D:/Personal Projects Git/DocKeep/Charts/4.PNG"

# Attach the 4-hour chart screenshot at column I / row 2, matching the
# same size used for every other chart picture in this workbook
# (10001250 x 5715000 EMU == 787.5 x 450 points).
$anchor = $ws.Cells.Item(2, 9)
$ws.Shapes.AddPicture("/tmp/work/extracted/xl/media/image6.png", $false, $true, $anchor.Left, $anchor.Top, 787.5, 450) | Out-Null
